# Sort the unit dictionary rows 12-26 alphabetically by the ENr code (column A).
# The header row (row 1) and rows 2-11 are already in correct alphabetical order
# and remain untouched. Rows 12-26 get re-sorted so that E_IDX2010 and
# E_IDX2015 move up (right after E_IDX2008) and E_NN moves between E_MNEW and
# E_NUM, with everything else shifting down accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("E_IDX2010",     "2010 = 100",                                                            "2010 = 100"),
    @("E_IDX2015",     "2015 = 100",                                                            "2015 = 100"),
    @("E_KGPHA",       "Kilogramm pro Hektar",                                                  "Kilogram per hectare"),
    @("E_MILIGPL",     "Miligramm pro Liter",                                                   "Miligrams per litre"),
    @("E_MIN",         "Minuten",                                                                "Minutes"),
    @("E_MNEUR",       "Millionen EUR",                                                          "Million EUR"),
    @("E_MNEW",        "Millionen Einwohner/ -innen",                                            "Million inhabitants"),
    @("E_NN",          "",                                                                       ""),
    @("E_NUM",         "Anzahl",                                                                 "Number"),
    @("E_P10H6EWN",    "Je 100 000 Einwohner/ -innen",                                           "Per 100,000 inhabitants"),
    @("E_P10H6EWNU70", "Je 100 000 Einwohner/ -innen unter 70 Jahren (ohne unter 1-Jährige)",    "Per 100,000 inhabitants under 70 years (excluding under 1 year olds)"),
    @("E_PRZNT",       "Prozent",                                                                "Percentage"),
    @("E_QMPA",        "m² pro Jahr",                                                            "m² per year"),
    @("E_QMPINHABA",   "m² pro Einwohner/-in und Jahr",                                          "m² per inhabitant and year"),
    @("E_TEUR",        "1 000 EUR",                                                              "1.000 EUR")
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
